$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$bvals1 = @(49,51,52,54,56,58,59,61,63,65,66,68,70,71,73,75,77,78,80,82,84,86,88,90,92,94,96,98,100,103,106,108,112,115,120,129,130,130,130,130,130)
for ($r = 0; $r -lt $bvals1.Length; $r++) {
    $ws.Cells.Item($r+2, 1).Value = $r
    $ws.Cells.Item($r+2, 2).Value = $bvals1[$r]
}

$ws = $wb.Worksheets.Item(2)
$bvals2 = @(49,50,52,54,55,57,59,60,62,64,66,67,69,71,72,74,76,77,79,81,83,85,86,88,90,92,94,96,99,101,103,106,109,112,116,120,128,130,130,130,130)
for ($r = 0; $r -lt $bvals2.Length; $r++) {
    $ws.Cells.Item($r+2, 1).Value = $r
    $ws.Cells.Item($r+2, 2).Value = $bvals2[$r]
}

$ws = $wb.Worksheets.Item(3)
$bvals3 = @(49,50,52,54,55,57,59,60,62,64,65,67,69,70,72,74,75,77,79,80,82,84,86,88,90,91,93,95,98,100,102,104,107,110,113,117,121,130,130,130,130)
for ($r = 0; $r -lt $bvals3.Length; $r++) {
    $ws.Cells.Item($r+2, 1).Value = $r
    $ws.Cells.Item($r+2, 2).Value = $bvals3[$r]
}

$ws = $wb.Worksheets.Item(4)
$bvals4 = @(49,51,53,54,56,58,59,61,63,64,66,67,69,71,72,74,76,77,79,81,82,84,86,88,90,91,93,95,97,100,102,104,107,109,112,116,120,126,130,130,130)
for ($r = 0; $r -lt $bvals4.Length; $r++) {
    $ws.Cells.Item($r+2, 1).Value = $r
    $ws.Cells.Item($r+2, 2).Value = $bvals4[$r]
}

$ws = $wb.Worksheets.Item(5)
$bvals5 = @(51,52,54,56,57,59,60,62,64,65,67,68,70,72,73,75,77,78,80,82,83,85,87,89,90,92,94,96,98,100,103,105,107,110,113,117,121,127,130,130,130)
for ($r = 0; $r -lt $bvals5.Length; $r++) {
    $ws.Cells.Item($r+2, 1).Value = $r
    $ws.Cells.Item($r+2, 2).Value = $bvals5[$r]
}

$ws = $wb.Worksheets.Item(6)
$bvals6 = @(54,56,57,59,60,62,63,65,67,68,70,71,73,74,76,78,79,81,83,84,86,88,89,91,93,95,97,99,101,103,106,108,111,114,117,121,127,130,130,130,130)
for ($r = 0; $r -lt $bvals6.Length; $r++) {
    $ws.Cells.Item($r+2, 1).Value = $r
    $ws.Cells.Item($r+2, 2).Value = $bvals6[$r]
}

Write-Host "done"